# Add season record columns (Wins / Losses / Ties) to the CLE_2023 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold, bordered, centered) from A1 onto the
# three new header cells so they reuse the same style as the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (same for every player row) for rows 2-52.
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 76   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 86   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
